$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 16
$ws.Range("K2").Value = 2.666666666666667
$ws.Range("L2").Value = $false
$ws.Range("P2").Value = $false
$ws.Range("Q2").Value = $true
$ws.Range("R2").Value = 0
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 1
$ws.Range("Y2").Value = 0

# Row 3
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 6
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 9
$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 35
$ws.Range("K3").Value = 5.833333333333333
$ws.Range("M3").Value = $false
$ws.Range("P3").Value = $true
$ws.Range("S3").Value = 0
$ws.Range("V3").Value = 1
$ws.Range("Y3").Value = 1

# Row 4
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = 6
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 7
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 29
$ws.Range("K4").Value = 4.833333333333333
$ws.Range("M4").Value = $true
$ws.Range("N4").Value = $false
$ws.Range("Q4").Value = $false
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("X4").Value = 2

# Row 5
$ws.Range("D5").Value = 7
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 7
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 21
$ws.Range("K5").Value = 3.5
$ws.Range("L5").Value = $true
$ws.Range("N5").Value = $true
$ws.Range("R5").Value = 1
$ws.Range("T5").Value = 1
$ws.Range("X5").Value = 2

# Row 6
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 7
$ws.Range("I6").Value = 6
$ws.Range("J6").Value = 25
$ws.Range("K6").Value = 4.166666666666667
$ws.Range("P6").Value = $false
$ws.Range("V6").Value = 0
$ws.Range("X6").Value = 0
